$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Values are written in this specific order so that new shared-string
# table entries are appended in the same order as the target workbook.
$ws.Range("A13").Value = "18/12/2019"
$ws.Range("C13").Value = "table design for project was done"

$ws.Range("A14").Value = "19/12/2019"

$ws.Range("B13").Value = "the screens for modules were started"

$ws.Range("B14").Value = "learned new concepts such as data strucure,technolgies,full stack developer,domain"

$ws.Range("A15").Value = "20/12/2019"

$ws.Range("B15").Value = "brush up of java concepts,jsp"
$ws.Range("C15").Value = "brush up of java concepts,jsp"

$ws.Range("A16").Value = "23/12/2019"

$ws.Range("B16").Value = "went through oops concepts such as polymorphism,encapsulation,classes etc……"

$ws.Range("B21").Select()
